# Catalogo de productos - add "Existencia Almacen" column (D) with stock data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 2: move the "Fecha:" label and the =TODAY() formula one column
# to the right (B2 -> C2, C2 -> D2) to make room for the new column.
# ------------------------------------------------------------------
$ws.Range("C2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D2").Formula = "=TODAY()"

$ws.Range("B2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C2").Value = "Fecha:"

$ws.Range("B2").Clear() | Out-Null

# ------------------------------------------------------------------
# Row 4: new header for the added column
# ------------------------------------------------------------------
$ws.Range("C4").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D4").Value = "Existencia Almacen"

# ------------------------------------------------------------------
# Body rows 5:40 - give column D the same plain style used across the
# rest of the table
# ------------------------------------------------------------------
$ws.Range("A5").Copy() | Out-Null
$ws.Range("D5:D40").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Rows 7 and 8 in column D pick up the bordered style that column C
# already uses on row 8
$ws.Range("C8").Copy() | Out-Null
$ws.Range("D7:D8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ------------------------------------------------------------------
# Row 1 / Row 3: touch column D so the sheet's used range/dimension
# grows to include it, without altering any visible formatting
# ------------------------------------------------------------------
$ws.Range("D1").Borders.LineStyle = -4142   # xlLineStyleNone

# ------------------------------------------------------------------
# Column widths
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 53
$ws.Columns.Item(3).ColumnWidth = 19
$ws.Columns.Item(4).ColumnWidth = 15

# ------------------------------------------------------------------
# Selection as left by the author
# ------------------------------------------------------------------
$ws.Range("D5").Select() | Out-Null

$wb.Save()
